$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1024.75
$ws.Range("I15").Value = 1024.75
$ws.Range("K15").Value = 3074.25
$ws.Range("M15").Value = -2905.25
$ws.Range("H17").Value = 1047.6111
$ws.Range("J17").Value = 858.6415
$ws.Range("L17").Value = 2575.9245
$ws.Range("N17").Value = -2911.9245
$ws.Range("H41").Value = 451.54544
$ws.Range("I41").Value = 373.66666
$ws.Range("K41").Value = 373.66666
$ws.Range("M41").Value = 66.33334000000002
$ws.Range("H62").Value = 2170.182
$ws.Range("I62").Value = 1485.75
$ws.Range("J62").Value = 3995.3333
$ws.Range("K62").Value = 1485.75
$ws.Range("L62").Value = 3995.3333
$ws.Range("M62").Value = -861.75
$ws.Range("N62").Value = -5243.3333
$ws.Range("H65").Value = 2170.182
$ws.Range("I65").Value = 1485.75
$ws.Range("J65").Value = 3995.3333
$ws.Range("K65").Value = 7428.75
$ws.Range("L65").Value = 19976.6665
$ws.Range("M65").Value = -4308.75
$ws.Range("N65").Value = -26216.6665
$ws.Range("H98").Value = 2227
$ws.Range("I98").Value = 2450.0908
$ws.Range("K98").Value = 2450.0908
$ws.Range("M98").Value = -952.0907999999999
$ws.Range("H113").Value = 15063.444
$ws.Range("I113").Value = 16508.875
$ws.Range("K113").Value = 16508.875
$ws.Range("M113").Value = -13254.875
$ws.Range("H122").Value = 2227
$ws.Range("I122").Value = 2450.0908
$ws.Range("K122").Value = 7350.2724
$ws.Range("M122").Value = -4900.2724
$ws.Range("H131").Value = 1509.5385
$ws.Range("I131").Value = 696.9
$ws.Range("J131").Value = 4218.3335
$ws.Range("K131").Value = 2090.7
$ws.Range("L131").Value = 12655.0005
$ws.Range("M131").Value = 2949.3
$ws.Range("N131").Value = -22735.0005
$ws.Range("H132").Value = 1060.8611
$ws.Range("I132").Value = 1070.3235
$ws.Range("K132").Value = 3210.9705
$ws.Range("M132").Value = -680.9704999999999
$ws.Range("H135").Value = 916
$ws.Range("I135").Value = 851.1111
$ws.Range("K135").Value = 7659.9999
$ws.Range("M135").Value = -5124.9999
$ws.Range("H137").Value = 1694.5
$ws.Range("I137").Value = 1212.75
$ws.Range("J137").Value = 1969.7858
$ws.Range("K137").Value = 3638.25
$ws.Range("L137").Value = 5909.357400000001
$ws.Range("M137").Value = -1088.25
$ws.Range("N137").Value = -11009.3574
$ws.Range("H141").Value = 3912.923
$ws.Range("I141").Value = 2648.1667
$ws.Range("K141").Value = 7944.500100000001
$ws.Range("M141").Value = -2764.500100000001

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3606.5425
$ws.Range("I32").Value = 2343.587
$ws.Range("J32").Value = 8075.4614
$ws.Range("K32").Value = 2343.587
$ws.Range("L32").Value = 8075.4614
$ws.Range("M32").Value = -2056.587
$ws.Range("N32").Value = -8649.4614
$ws.Range("H45").Value = 1654.5264
$ws.Range("I45").Value = 840.3333
$ws.Range("K45").Value = 840.3333
$ws.Range("M45").Value = -463.3333
$ws.Range("H61").Value = 3116.9333
$ws.Range("I61").Value = 1617.4166
$ws.Range("K61").Value = 1617.4166
$ws.Range("M61").Value = -1405.4166
$ws.Range("H74").Value = 500.57144
$ws.Range("I74").Value = 500.57144
$ws.Range("K74").Value = 500.57144
$ws.Range("M74").Value = 373.42856
$ws.Range("H77").Value = 500.57144
$ws.Range("I77").Value = 500.57144
$ws.Range("K77").Value = 2502.8572
$ws.Range("M77").Value = 1865.1428
$ws.Range("H114").Value = 0
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H132").Value = 2450.7
$ws.Range("I132").Value = 1813.625
$ws.Range("K132").Value = 5440.875
$ws.Range("M132").Value = -2910.875
$ws.Range("H135").Value = 50000
$ws.Range("J135").Value = 50000
$ws.Range("L135").Value = 50000
$ws.Range("N135").Value = -60140
$ws.Range("H136").Value = 3116.9333
$ws.Range("I136").Value = 1617.4166
$ws.Range("K136").Value = 4852.2498
$ws.Range("M136").Value = -2302.2498

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1872.4445
$ws.Range("I20").Value = 1651.0714
$ws.Range("K20").Value = 1651.0714
$ws.Range("M20").Value = -1404.0714
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("H134").Value = 7667.5527
$ws.Range("I134").Value = 8852
$ws.Range("K134").Value = 26556
$ws.Range("M134").Value = -24021

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 633.3333
$ws.Range("I19").Value = 633.3333
$ws.Range("K19").Value = 633.3333
$ws.Range("M19").Value = -463.3333
$ws.Range("H24").Value = 633.3333
$ws.Range("I24").Value = 633.3333
$ws.Range("K24").Value = 633.3333
$ws.Range("M24").Value = -463.3333
$ws.Range("H51").Value = 30085.715
$ws.Range("J51").Value = 30085.715
$ws.Range("L51").Value = 30085.715
$ws.Range("N51").Value = -31557.715
$ws.Range("H61").Value = 30085.715
$ws.Range("J61").Value = 30085.715
$ws.Range("L61").Value = 30085.715
$ws.Range("N61").Value = -30781.715
$ws.Range("H99").Value = 2268.75
$ws.Range("H126").Value = 2268.75
$ws.Range("H132").Value = 2157.04
$ws.Range("I132").Value = 1136.1052
$ws.Range("K132").Value = 3408.3156
$ws.Range("M132").Value = -878.3155999999999
$ws.Range("H134").Value = 715.8095

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 328.57144
$ws.Range("J92").Value = 328.57144
$ws.Range("L92").Value = 985.71432
$ws.Range("N92").Value = -3481.71432
$ws.Range("H97").Value = 841.5
$ws.Range("J97").Value = 847
$ws.Range("L97").Value = 2541
$ws.Range("N97").Value = -3533
$ws.Range("H131").Value = 10538.329
$ws.Range("J131").Value = 11269.735
$ws.Range("L131").Value = 33809.205
$ws.Range("N131").Value = -43889.205
$ws.Range("H132").Value = 1511.1111
$ws.Range("I132").Value = 1500
$ws.Range("K132").Value = 13500
$ws.Range("M132").Value = -10970
$ws.Range("H140").Value = 2090.2727
$ws.Range("I140").Value = 1038
$ws.Range("K140").Value = 3114
$ws.Range("M140").Value = 2066

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1981.25
$ws.Range("I122").Value = 1731.8948
$ws.Range("J122").Value = 2345.6924
$ws.Range("K122").Value = 5195.6844
$ws.Range("L122").Value = 7037.0772
$ws.Range("M122").Value = -2745.6844
$ws.Range("N122").Value = -11937.0772
$ws.Range("H123").Value = 14260.6
$ws.Range("J123").Value = 14260.6
$ws.Range("L123").Value = 14260.6
$ws.Range("N123").Value = -19160.6

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7290.1
$ws.Range("I7").Value = 4668
$ws.Range("K7").Value = 4668
$ws.Range("M7").Value = -4556
$ws.Range("H46").Value = 1677.7778
$ws.Range("J46").Value = 1737.5
$ws.Range("L46").Value = 1737.5
$ws.Range("N46").Value = -2113.5
$ws.Range("H100").Value = 1800
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 1800
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 1800
$ws.Range("N100").Value = -2882
$ws.Range("M100").ClearContents()
$ws.Range("H126").Value = 7290.1
$ws.Range("I126").Value = 4668
$ws.Range("K126").Value = 14004
$ws.Range("M126").Value = -11534
$ws.Range("H132").Value = 1944.8695
$ws.Range("I132").Value = 1099
$ws.Range("K132").Value = 3297
$ws.Range("M132").Value = -767

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1111.7142
$ws.Range("I113").Value = 1025.5
$ws.Range("J113").Value = 1226.6666
$ws.Range("K113").Value = 3076.5
$ws.Range("L113").Value = 3679.9998
$ws.Range("M113").Value = -906.5
$ws.Range("N113").Value = -8019.9998
$ws.Range("H126").Value = 6090.6
$ws.Range("J126").Value = 7054.6665
$ws.Range("L126").Value = 21163.9995
$ws.Range("N126").Value = -26103.9995
$ws.Range("H132").Value = 4035.7144
$ws.Range("I132").Value = 3773
$ws.Range("K132").Value = 11319
$ws.Range("M132").Value = -8789

